# Update automàtic: dades i banners [2026-02-09 20:50]
#
# Refreshes the meteocat daily-summary scrape: each station row gets a new
# DATA_EXTRACCIO timestamp plus whichever measurement columns the upstream
# service revised (precipitation, pressure, humidity, snow depth, mean temp).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal-text value (e.g. "86%") without letting Excel
# reinterpret it as a percentage number and reformat the cell. We build the
# text via a formula in an unused scratch cell (so it is a genuine string,
# never subject to "smart" numeric-entry parsing), then copy / paste-special
# just the value into the destination, leaving its existing style untouched.
function Set-LiteralText($cellRef, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# --- Plain text / timestamp updates -------------------------------------
$ws.Range('E2').Value = '2026-02-09 20:48:32'
$ws.Range('I2').Value = '1.4 mm'
$ws.Range('O2').Value = '-0.7 °C'
$ws.Range('E3').Value = '2026-02-09 20:48:35'
$ws.Range('I3').Value = '2.1 mm'
$ws.Range('O3').Value = '-3.4 °C'
$ws.Range('E4').Value = '2026-02-09 20:48:37'
$ws.Range('J4').Value = '1007.0 hPa'
$ws.Range('E5').Value = '2026-02-09 20:48:40'
$ws.Range('O5').Value = '-2.6 °C'
$ws.Range('E6').Value = '2026-02-09 20:48:43'
$ws.Range('J6').Value = '1006.9 hPa'
$ws.Range('E7').Value = '2026-02-09 20:48:45'
$ws.Range('E8').Value = '2026-02-09 20:48:48'
$ws.Range('E9').Value = '2026-02-09 20:48:51'
$ws.Range('E10').Value = '2026-02-09 20:48:54'
$ws.Range('O10').Value = '8.1 °C'
$ws.Range('E11').Value = '2026-02-09 20:48:56'
$ws.Range('O11').Value = '5.3 °C'
$ws.Range('E12').Value = '2026-02-09 20:48:59'
$ws.Range('O12').Value = '8.9 °C'
$ws.Range('E13').Value = '2026-02-09 20:49:01'
$ws.Range('G13').Value = '1 cm'
$ws.Range('E14').Value = '2026-02-09 20:49:04'
$ws.Range('E15').Value = '2026-02-09 20:49:06'
$ws.Range('O15').Value = '8.2 °C'
$ws.Range('E16').Value = '2026-02-09 20:49:09'
$ws.Range('I16').Value = '1.6 mm'
$ws.Range('E17').Value = '2026-02-09 20:49:11'
$ws.Range('E18').Value = '2026-02-09 20:49:14'
$ws.Range('O18').Value = '9.0 °C'
$ws.Range('E19').Value = '2026-02-09 20:49:17'
$ws.Range('E20').Value = '2026-02-09 20:49:19'
$ws.Range('I20').Value = '0.1 mm'
$ws.Range('O20').Value = '-4.3 °C'
$ws.Range('E21').Value = '2026-02-09 20:49:22'
$ws.Range('O21').Value = '4.7 °C'
$ws.Range('E22').Value = '2026-02-09 20:49:24'
$ws.Range('O22').Value = '-5.0 °C'
$ws.Range('E23').Value = '2026-02-09 20:49:27'
$ws.Range('I23').Value = '0.9 mm'
$ws.Range('E24').Value = '2026-02-09 20:49:29'
$ws.Range('I24').Value = '1.2 mm'
$ws.Range('E25').Value = '2026-02-09 20:49:31'
$ws.Range('E26').Value = '2026-02-09 20:49:34'
$ws.Range('E27').Value = '2026-02-09 20:49:37'
$ws.Range('I27').Value = '0.1 mm'
$ws.Range('E28').Value = '2026-02-09 20:49:39'
$ws.Range('E29').Value = '2026-02-09 20:49:42'
$ws.Range('E30').Value = '2026-02-09 20:49:45'
$ws.Range('E31').Value = '2026-02-09 20:49:47'
$ws.Range('E32').Value = '2026-02-09 20:49:50'
$ws.Range('I32').Value = '0.7 mm'
$ws.Range('E33').Value = '2026-02-09 20:49:52'
$ws.Range('E34').Value = '2026-02-09 20:49:55'
$ws.Range('E35').Value = '2026-02-09 20:49:58'
$ws.Range('I35').Value = '1.6 mm'
$ws.Range('E36').Value = '2026-02-09 20:50:00'
$ws.Range('J36').Value = '1007.1 hPa'
$ws.Range('E37').Value = '2026-02-09 20:50:03'
$ws.Range('E38').Value = '2026-02-09 20:50:06'
$ws.Range('O38').Value = '8.6 °C'
$ws.Range('E39').Value = '2026-02-09 20:50:08'
$ws.Range('E40').Value = '2026-02-09 20:50:11'
$ws.Range('O40').Value = '4.8 °C'
$ws.Range('E41').Value = '2026-02-09 20:50:14'
$ws.Range('I41').Value = '0.1 mm'
$ws.Range('J41').Value = '1007.7 hPa'
$ws.Range('E42').Value = '2026-02-09 20:50:16'
$ws.Range('O42').Value = '9.0 °C'
$ws.Range('E43').Value = '2026-02-09 20:50:19'
$ws.Range('E44').Value = '2026-02-09 20:50:21'
$ws.Range('I44').Value = '1.0 mm'
$ws.Range('E45').Value = '2026-02-09 20:50:23'
$ws.Range('I45').Value = '0.5 mm'
$ws.Range('O45').Value = '4.1 °C'
$ws.Range('E46').Value = '2026-02-09 20:50:26'

# --- Percentage-looking text updates (kept as literal text) -------------
Set-LiteralText 'H12' '86%'
Set-LiteralText 'H16' '73%'
Set-LiteralText 'H24' '83%'
Set-LiteralText 'H27' '83%'
Set-LiteralText 'H29' '84%'
Set-LiteralText 'H30' '86%'
Set-LiteralText 'H35' '78%'
Set-LiteralText 'H37' '80%'
Set-LiteralText 'H43' '73%'

